$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated cryptocurrency price/volume data scraped on Wed May 22 13:45:50 UTC 2024
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.619.51'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.682.70'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '614.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.59'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.672.48'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.46%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.529'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.10%  '
$ws.Range('E10').Value = '  -4.20%  '
$ws.Range('E11').Value = '  -2.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.478'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.26%  '
$ws.Range('E13').Value = '  -3.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000252'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.299.54'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.684.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.624.48'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.121'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.52'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '499.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.11%  '
$ws.Range('E23').Value = '  -5.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.50'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.33'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000128'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.07%  '
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.42'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.87'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.93'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.99'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.113'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.01'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.136'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.336'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.24%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.05'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.77%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '49.83'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.70%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '429.11'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.91'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.95'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.53'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.931.81'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0357'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.86%  '
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.70'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.46%  '
$ws.Range('E51').Value = '  -2.89%  '
